$wb = $excel.ActiveWorkbook

# --- "Links" sheet: clear the first-name value for Sriram (Sri) Panyam (row 35) ---
$wsLinks = $wb.Worksheets.Item("Links")
$wsLinks.Range("D35").Value = ""

# --- "Usuários" sheet: update the account-status info for row 13 / row 14 ---
$wsUsuarios = $wb.Worksheets.Item("Usuários")

# Row 13 (deodoro_fonseca@gmail.com): attempt count bumped to 10, account now blocked
$wsUsuarios.Range("D13").Value = 10
$wsUsuarios.Range("E13").Value = "Não"
$wsUsuarios.Range("F13").Value = "Não"
$wsUsuarios.Range("G13").Value = "Conta bloqueada pelo Linkedin por muitas tentativas. Troque esta conta por outra, ou remova esta linha do Excel."

# Row 14 (pedro_alvarez@gmail.com): attempt count bumped to 9, account now blocked
$wsUsuarios.Range("D14").Value = 9
$wsUsuarios.Range("E14").Value = "Não"
$wsUsuarios.Range("F14").Value = "Não"
$wsUsuarios.Range("G14").Value = "Conta bloqueada pelo Linkedin por muitas tentativas. Troque esta conta por outra, ou remova esta linha do Excel."

$wb.Save()
